$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

$newRows = @(
    @("vaishali.kh2310@gmail.com", "Login", "2025-06-16 23:00:07"),
    @("vaishali.kh2310@gmail.com", "Login", "2025-06-17 08:40:54"),
    @("6377384840vk@gmail.com",    "Login", "2025-06-17 11:35:12"),
    @("vaishali.kh2310@gmail.com", "Login", "2025-06-17 17:10:13")
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}
